$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue $ws "D2" "304.94"
Set-TextValue $ws "E2" "1.35%"
Set-TextValue $ws "E3" "-3.60%"
Set-TextValue $ws "D4" "5.038"
Set-TextValue $ws "E4" "1.31%"
Set-TextValue $ws "D5" "0.07843"
Set-TextValue $ws "E5" "1.36%"
Set-TextValue $ws "D6" "2.165"
Set-TextValue $ws "E6" "-1.71%"
Set-TextValue $ws "D7" "7.909"
Set-TextValue $ws "E7" "-1.14%"
Set-TextValue $ws "D8" "4.084"
Set-TextValue $ws "E8" "2.37%"
Set-TextValue $ws "D9" "0.9189"
Set-TextValue $ws "E9" "0.65%"
Set-TextValue $ws "D10" "0.09730"
Set-TextValue $ws "E10" "4.00%"
Set-TextValue $ws "D11" "0.1862"
Set-TextValue $ws "E11" "3.72%"
Set-TextValue $ws "D12" "0.08710"
Set-TextValue $ws "E12" "3.58%"
Set-TextValue $ws "D13" "0.03478"
Set-TextValue $ws "E13" "-1.59%"
Set-TextValue $ws "D14" "0.09906"
Set-TextValue $ws "E14" "-0.18%"
Set-TextValue $ws "D15" "0.001426"
Set-TextValue $ws "E15" "-2.76%"
Set-TextValue $ws "D16" "0.005716"
Set-TextValue $ws "E16" "0.45%"
Set-TextValue $ws "D17" "3.463"
Set-TextValue $ws "E17" "-0.36%"
Set-TextValue $ws "D18" "2.392"
Set-TextValue $ws "E18" "16.54%"
Set-TextValue $ws "E19" "-1.13%"
Set-TextValue $ws "E20" "2.33%"
Set-TextValue $ws "D21" "4.760"
Set-TextValue $ws "E21" "4.08%"
Set-TextValue $ws "D22" "0.2204"
Set-TextValue $ws "E22" "-0.97%"
Set-TextValue $ws "D23" "0.04601"
Set-TextValue $ws "E23" "-1.00%"
Set-TextValue $ws "D24" "0.005096"
Set-TextValue $ws "E24" "14.78%"
Set-TextValue $ws "D25" "0.001227"
Set-TextValue $ws "E26" "7.79%"
Set-TextValue $ws "E27" "0.13%"
Set-TextValue $ws "D39" "0.01831"
Set-TextValue $ws "E39" "4.36%"
Set-TextValue $ws "D40" "0.04764"
Set-TextValue $ws "E40" "1.84%"
Set-TextValue $ws "D41" "0.007681"
Set-TextValue $ws "E41" "-1.52%"
Set-TextValue $ws "D42" "0.1396"
Set-TextValue $ws "E42" "0.61%"
Set-TextValue $ws "D43" "0.007730"
Set-TextValue $ws "E43" "1.01%"
Set-TextValue $ws "E44" "-2.54%"
Set-TextValue $ws "D45" "0.01113"
Set-TextValue $ws "E45" "10.01%"
Set-TextValue $ws "D46" "0.00006367"
Set-TextValue $ws "E46" "4.63%"
Set-TextValue $ws "D48" "0.0005796"
Set-TextValue $ws "E48" "-0.08%"
Set-TextValue $ws "D49" "24.33"
Set-TextValue $ws "E49" "180.96%"
